$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H33").Value = 128.14285
$ws_ALC.Range("I33").Value = 130.3077
$ws_ALC.Range("J33").Value = 100
$ws_ALC.Range("K33").Value = 130.3077
$ws_ALC.Range("L33").Value = 100
$ws_ALC.Range("M33").Value = 98.69229999999999
$ws_ALC.Range("N33").Value = -558

$ws_ALC.Range("H38").Value = 3806.6
$ws_ALC.Range("I38").Value = 521.5
$ws_ALC.Range("J38").Value = 5996.6665
$ws_ALC.Range("K38").Value = 1564.5
$ws_ALC.Range("L38").Value = 17989.9995
$ws_ALC.Range("M38").Value = -1192.5
$ws_ALC.Range("N38").Value = -18733.9995

$ws_ALC.Range("H80").Value = 329.44
$ws_ALC.Range("I80").Value = 308.125
$ws_ALC.Range("J80").Value = 367.33334
$ws_ALC.Range("K80").Value = 924.375
$ws_ALC.Range("L80").Value = 1102.00002
$ws_ALC.Range("M80").Value = 73.625
$ws_ALC.Range("N80").Value = -3098.00002

$ws_ALC.Range("H83").Value = 329.44
$ws_ALC.Range("I83").Value = 308.125
$ws_ALC.Range("J83").Value = 367.33334
$ws_ALC.Range("K83").Value = 2773.125
$ws_ALC.Range("L83").Value = 3306.00006
$ws_ALC.Range("M83").Value = 2218.875
$ws_ALC.Range("N83").Value = -13290.00006

$ws_ALC.Range("H86").Value = 5998.4
$ws_ALC.Range("I86").Value = 5996
$ws_ALC.Range("J86").Value = 5999
$ws_ALC.Range("K86").Value = 5996
$ws_ALC.Range("L86").Value = 5999
$ws_ALC.Range("M86").Value = -4873
$ws_ALC.Range("N86").Value = -8245

$ws_ALC.Range("H89").Value = 5998.4
$ws_ALC.Range("I89").Value = 5996
$ws_ALC.Range("J89").Value = 5999
$ws_ALC.Range("K89").Value = 29980
$ws_ALC.Range("L89").Value = 29995
$ws_ALC.Range("M89").Value = -24364
$ws_ALC.Range("N89").Value = -41227

$ws_ALC.Range("H100").Value = 1145.3158
$ws_ALC.Range("I100").Value = 976.5454999999999
$ws_ALC.Range("J100").Value = 1377.375
$ws_ALC.Range("K100").Value = 976.5454999999999
$ws_ALC.Range("L100").Value = 1377.375
$ws_ALC.Range("M100").Value = -435.5454999999999
$ws_ALC.Range("N100").Value = -2459.375

$ws_ALC.Range("H132").Value = 2775.5
$ws_ALC.Range("I132").Value = 2750.5557
$ws_ALC.Range("J132").Value = 3000
$ws_ALC.Range("K132").Value = 8251.667099999999
$ws_ALC.Range("L132").Value = 9000
$ws_ALC.Range("M132").Value = -5721.667099999999
$ws_ALC.Range("N132").Value = -14060

$ws_ALC.Range("H137").Value = 2775
$ws_ALC.Range("I137").Value = 1321.7333
$ws_ALC.Range("J137").Value = 4756.727
$ws_ALC.Range("K137").Value = 3965.199900000001
$ws_ALC.Range("L137").Value = 14270.181
$ws_ALC.Range("M137").Value = -1415.199900000001
$ws_ALC.Range("N137").Value = -19370.181

$ws_ARM.Range("H2").Value = 4995
$ws_ARM.Range("I2").Value = 4995
$ws_ARM.Range("J2").Value = 0
$ws_ARM.Range("K2").Value = 4995
$ws_ARM.Range("L2").Value = 0
$ws_ARM.Range("M2").Value = $null
$ws_ARM.Range("N2").Value = -4882

$ws_ARM.Range("H19").Value = 386.5
$ws_ARM.Range("I19").Value = 386.5
$ws_ARM.Range("J19").Value = 0
$ws_ARM.Range("K19").Value = 386.5
$ws_ARM.Range("L19").Value = 0
$ws_ARM.Range("M19").Value = $null
$ws_ARM.Range("N19").Value = -157.5

$ws_ARM.Range("H30").Value = 10269.667
$ws_ARM.Range("I30").Value = 15199.5
$ws_ARM.Range("J30").Value = 410
$ws_ARM.Range("K30").Value = 15199.5
$ws_ARM.Range("L30").Value = 410
$ws_ARM.Range("M30").Value = -15049.5
$ws_ARM.Range("N30").Value = -710

$ws_ARM.Range("H76").Value = 17500
$ws_ARM.Range("I76").Value = 0
$ws_ARM.Range("J76").Value = 17500
$ws_ARM.Range("K76").Value = 0
$ws_ARM.Range("L76").Value = 17500
$ws_ARM.Range("N76").Value = -18176

$ws_ARM.Range("H79").Value = 17500
$ws_ARM.Range("I79").Value = 0
$ws_ARM.Range("J79").Value = 17500
$ws_ARM.Range("K79").Value = 0
$ws_ARM.Range("L79").Value = 17500
$ws_ARM.Range("N79").Value = -19840

$ws_ARM.Range("H110").Value = 17395.6
$ws_ARM.Range("I110").Value = 19422.666
$ws_ARM.Range("J110").Value = 14355
$ws_ARM.Range("K110").Value = 19422.666
$ws_ARM.Range("L110").Value = 14355
$ws_ARM.Range("M110").Value = -17377.666
$ws_ARM.Range("N110").Value = -18445

$ws_ARM.Range("H116").Value = 4995
$ws_ARM.Range("I116").Value = 4995
$ws_ARM.Range("J116").Value = 0
$ws_ARM.Range("K116").Value = 4995
$ws_ARM.Range("L116").Value = 0
$ws_ARM.Range("M116").Value = $null
$ws_ARM.Range("N116").Value = -2701

$ws_ARM.Range("H132").Value = 1402.7307
$ws_ARM.Range("I132").Value = 1402.7307
$ws_ARM.Range("J132").Value = 0
$ws_ARM.Range("K132").Value = 4208.1921
$ws_ARM.Range("L132").Value = 0
$ws_ARM.Range("M132").Value = -1678.1921

$ws_BSM.Range("H3").Value = 4995
$ws_BSM.Range("I3").Value = 4995
$ws_BSM.Range("J3").Value = 0
$ws_BSM.Range("K3").Value = 4995
$ws_BSM.Range("L3").Value = 0
$ws_BSM.Range("M3").Value = $null
$ws_BSM.Range("N3").Value = -4881

$ws_BSM.Range("H22").Value = 699.75
$ws_BSM.Range("I22").Value = 699.75
$ws_BSM.Range("J22").Value = 0
$ws_BSM.Range("K22").Value = 699.75
$ws_BSM.Range("L22").Value = 0
$ws_BSM.Range("M22").Value = -526.75

$ws_BSM.Range("H80").Value = 735.8
$ws_BSM.Range("I80").Value = 953.3333
$ws_BSM.Range("J80").Value = 409.5
$ws_BSM.Range("K80").Value = 953.3333
$ws_BSM.Range("L80").Value = 409.5
$ws_BSM.Range("M80").Value = 44.66669999999999
$ws_BSM.Range("N80").Value = -2405.5

$ws_BSM.Range("H81").Value = 40000
$ws_BSM.Range("I81").Value = 40000
$ws_BSM.Range("J81").Value = 0
$ws_BSM.Range("K81").Value = 40000
$ws_BSM.Range("L81").Value = 0
$ws_BSM.Range("M81").Value = -38939
$ws_BSM.Range("N81").Value = $null

$ws_BSM.Range("H83").Value = 735.8
$ws_BSM.Range("I83").Value = 953.3333
$ws_BSM.Range("J83").Value = 409.5
$ws_BSM.Range("K83").Value = 4766.6665
$ws_BSM.Range("L83").Value = 2047.5
$ws_BSM.Range("M83").Value = 225.3334999999997
$ws_BSM.Range("N83").Value = -12031.5

$ws_BSM.Range("H84").Value = 40000
$ws_BSM.Range("I84").Value = 40000
$ws_BSM.Range("J84").Value = 0
$ws_BSM.Range("K84").Value = 120000
$ws_BSM.Range("L84").Value = 0
$ws_BSM.Range("M84").Value = -114696
$ws_BSM.Range("N84").Value = $null

$ws_BSM.Range("H86").Value = 2172.4443
$ws_BSM.Range("I86").Value = 1975.6666
$ws_BSM.Range("J86").Value = 2566
$ws_BSM.Range("K86").Value = 1975.6666
$ws_BSM.Range("L86").Value = 2566
$ws_BSM.Range("M86").Value = -852.6666
$ws_BSM.Range("N86").Value = -4812

$ws_BSM.Range("H89").Value = 2172.4443
$ws_BSM.Range("I89").Value = 1975.6666
$ws_BSM.Range("J89").Value = 2566
$ws_BSM.Range("K89").Value = 9878.333000000001
$ws_BSM.Range("L89").Value = 12830
$ws_BSM.Range("M89").Value = -4262.333000000001
$ws_BSM.Range("N89").Value = -24062

$ws_BSM.Range("H93").Value = 26000
$ws_BSM.Range("I93").Value = 0
$ws_BSM.Range("J93").Value = 26000
$ws_BSM.Range("K93").Value = 0
$ws_BSM.Range("L93").Value = 26000
$ws_BSM.Range("N93").Value = -29744

$ws_BSM.Range("H105").Value = 4311.857
$ws_BSM.Range("I105").Value = 3027.5
$ws_BSM.Range("J105").Value = 29999
$ws_BSM.Range("K105").Value = 3027.5
$ws_BSM.Range("L105").Value = 29999
$ws_BSM.Range("M105").Value = -1280.5
$ws_BSM.Range("N105").Value = -33493

$ws_BSM.Range("H130").Value = 63525.715
$ws_BSM.Range("I130").Value = 20000
$ws_BSM.Range("J130").Value = 70780
$ws_BSM.Range("K130").Value = 20000
$ws_BSM.Range("L130").Value = 70780
$ws_BSM.Range("M130").Value = -14980
$ws_BSM.Range("N130").Value = -80820

$ws_CRP.Range("H7").Value = 222.76923
$ws_CRP.Range("I7").Value = 199.66667
$ws_CRP.Range("J7").Value = 500
$ws_CRP.Range("K7").Value = 199.66667
$ws_CRP.Range("L7").Value = 500
$ws_CRP.Range("M7").Value = -86.66667000000001
$ws_CRP.Range("N7").Value = -726

$ws_CRP.Range("H22").Value = 685.5
$ws_CRP.Range("I22").Value = 446
$ws_CRP.Range("J22").Value = 765.3333
$ws_CRP.Range("K22").Value = 446
$ws_CRP.Range("L22").Value = 765.3333
$ws_CRP.Range("M22").Value = -96
$ws_CRP.Range("N22").Value = -1465.3333

$ws_CRP.Range("H122").Value = 3150.2104
$ws_CRP.Range("I122").Value = 3255.4
$ws_CRP.Range("J122").Value = 2755.75
$ws_CRP.Range("K122").Value = 9766.200000000001
$ws_CRP.Range("L122").Value = 8267.25
$ws_CRP.Range("M122").Value = -7316.200000000001
$ws_CRP.Range("N122").Value = -13167.25

$ws_CUL.Range("H44").Value = 1524.3
$ws_CUL.Range("I44").Value = 710.75
$ws_CUL.Range("J44").Value = 2066.6667
$ws_CUL.Range("K44").Value = 2132.25
$ws_CUL.Range("L44").Value = 6200.000100000001
$ws_CUL.Range("M44").Value = -1734.25
$ws_CUL.Range("N44").Value = -6996.000100000001

$ws_CUL.Range("H46").Value = 2650
$ws_CUL.Range("I46").Value = 300
$ws_CUL.Range("J46").Value = 5000
$ws_CUL.Range("K46").Value = 900
$ws_CUL.Range("L46").Value = 15000
$ws_CUL.Range("M46").Value = -809
$ws_CUL.Range("N46").Value = -15182

$ws_CUL.Range("H58").Value = 332.5
$ws_CUL.Range("I58").Value = 379.2
$ws_CUL.Range("J58").Value = 99
$ws_CUL.Range("K58").Value = 1137.6
$ws_CUL.Range("L58").Value = 297
$ws_CUL.Range("M58").Value = -1009.6
$ws_CUL.Range("N58").Value = -553

$ws_CUL.Range("H99").Value = 0
$ws_CUL.Range("I99").Value = 0
$ws_CUL.Range("J99").Value = 0
$ws_CUL.Range("K99").Value = 0
$ws_CUL.Range("L99").Value = $null
$ws_CUL.Range("N99").Value = 0

$ws_CUL.Range("H117").Value = 1417.3077
$ws_CUL.Range("I117").Value = 272.77777
$ws_CUL.Range("J117").Value = 3992.5
$ws_CUL.Range("K117").Value = 818.33331
$ws_CUL.Range("L117").Value = 11977.5
$ws_CUL.Range("M117").Value = 2623.66669
$ws_CUL.Range("N117").Value = -18861.5

$ws_CUL.Range("H129").Value = 581
$ws_CUL.Range("I129").Value = 298.2
$ws_CUL.Range("J129").Value = 1995
$ws_CUL.Range("K129").Value = 894.5999999999999
$ws_CUL.Range("L129").Value = 5985
$ws_CUL.Range("M129").Value = 4105.4
$ws_CUL.Range("N129").Value = -15985

$ws_CUL.Range("H139").Value = 2214.6667
$ws_CUL.Range("I139").Value = 1757.8
$ws_CUL.Range("J139").Value = 4499
$ws_CUL.Range("K139").Value = 5273.4
$ws_CUL.Range("L139").Value = 13497
$ws_CUL.Range("M139").Value = -133.3999999999996
$ws_CUL.Range("N139").Value = -23777

$ws_CUL.Range("H140").Value = 2389.6667
$ws_CUL.Range("I140").Value = 2389.6667
$ws_CUL.Range("J140").Value = 0
$ws_CUL.Range("K140").Value = 7169.000100000001
$ws_CUL.Range("L140").Value = 0
$ws_CUL.Range("M140").Value = -1989.000100000001

$ws_GSM.Range("H2").Value = 341.84616
$ws_GSM.Range("I2").Value = 85.42856999999999
$ws_GSM.Range("J2").Value = 641
$ws_GSM.Range("K2").Value = 85.42856999999999
$ws_GSM.Range("L2").Value = 641
$ws_GSM.Range("M2").Value = 27.57143000000001
$ws_GSM.Range("N2").Value = -867

$ws_GSM.Range("H57").Value = 29999
$ws_GSM.Range("I57").Value = 0
$ws_GSM.Range("J57").Value = 29999
$ws_GSM.Range("K57").Value = 0
$ws_GSM.Range("L57").Value = 29999
$ws_GSM.Range("N57").Value = -31639

$ws_GSM.Range("H80").Value = 3888.7778
$ws_GSM.Range("I80").Value = 4166.6665
$ws_GSM.Range("J80").Value = 3333
$ws_GSM.Range("K80").Value = 4166.6665
$ws_GSM.Range("L80").Value = 3333
$ws_GSM.Range("M80").Value = -3168.6665
$ws_GSM.Range("N80").Value = -5329

$ws_GSM.Range("H83").Value = 3888.7778
$ws_GSM.Range("I83").Value = 4166.6665
$ws_GSM.Range("J83").Value = 3333
$ws_GSM.Range("K83").Value = 20833.3325
$ws_GSM.Range("L83").Value = 16665
$ws_GSM.Range("M83").Value = -15841.3325
$ws_GSM.Range("N83").Value = -26649

$ws_GSM.Range("H97").Value = 1100.25
$ws_GSM.Range("I97").Value = 1068.5
$ws_GSM.Range("J97").Value = 1132
$ws_GSM.Range("K97").Value = 1068.5
$ws_GSM.Range("L97").Value = 1132
$ws_GSM.Range("M97").Value = -572.5
$ws_GSM.Range("N97").Value = -2124

$ws_GSM.Range("H113").Value = 5092.231
$ws_GSM.Range("I113").Value = 5097
$ws_GSM.Range("J113").Value = 5090.1113
$ws_GSM.Range("K113").Value = 5097
$ws_GSM.Range("L113").Value = 5090.1113
$ws_GSM.Range("M113").Value = -2927
$ws_GSM.Range("N113").Value = -9430.1113

$ws_LTW.Range("H46").Value = 2539.2
$ws_LTW.Range("I46").Value = 2197.8
$ws_LTW.Range("J46").Value = 2709.9
$ws_LTW.Range("K46").Value = 2197.8
$ws_LTW.Range("L46").Value = 2709.9
$ws_LTW.Range("M46").Value = -2009.8
$ws_LTW.Range("N46").Value = -3085.9

$ws_LTW.Range("H132").Value = 2399
$ws_LTW.Range("I132").Value = 1665
$ws_LTW.Range("J132").Value = 3500
$ws_LTW.Range("K132").Value = 4995
$ws_LTW.Range("L132").Value = 10500
$ws_LTW.Range("M132").Value = -2465
$ws_LTW.Range("N132").Value = -15560

$ws_WVR.Range("H60").Value = 69495
$ws_WVR.Range("I60").Value = 99994
$ws_WVR.Range("J60").Value = 38996
$ws_WVR.Range("K60").Value = 99994
$ws_WVR.Range("L60").Value = 38996
$ws_WVR.Range("M60").Value = -99172
$ws_WVR.Range("N60").Value = -40640
